$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new data row (row 13) for year "2021年", following the same
# layout as the existing rows (A = year label, B..U = indicator values,
# with columns G, Q and U left blank just like in the preceding rows).
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 214.6
$ws.Range("C13").Value = 273.9
$ws.Range("D13").Value = 237.5
$ws.Range("E13").Value = 269.8
$ws.Range("F13").Value = 260
$ws.Range("H13").Value = 175.8
$ws.Range("I13").Value = 162.6
$ws.Range("J13").Value = 272.3
$ws.Range("K13").Value = 175.6
$ws.Range("L13").Value = 252.9
$ws.Range("M13").Value = 172.9
$ws.Range("N13").Value = 211
$ws.Range("O13").Value = 254.4
$ws.Range("P13").Value = 216.2
$ws.Range("R13").Value = 157.1
$ws.Range("S13").Value = 273.5
$ws.Range("T13").Value = 145.9

# Touch the otherwise-blank cells so they persist as real (empty) cells in
# the sheet, matching the shape of the preceding rows (G/Q/U are blank
# there too).
$ws.Range("G13").Style = "Normal"
$ws.Range("Q13").Style = "Normal"
$ws.Range("U13").Style = "Normal"

# Give the new year label (column A) the same bold / centered / bordered
# look as the other year cells, by copying the formatting from A12 instead
# of re-building it (this reuses the existing style rather than creating a
# new, slightly different one).
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
